# "Kleine veranderingen in lay-out" - small layout changes across the 3 sheets.

$wb = $excel.ActiveWorkbook

$wsTester = $wb.Worksheets.Item("Tester")
$wsGen1   = $wb.Worksheets.Item("Generator1")
$wsGen2   = $wb.Worksheets.Item("Generator2")

# --- New shared string / label swap on Generator2 row 20 -------------------
# "4 channels, 5 comparators: 5346 networks" -> "...2376 networks"
$wsGen2.Range("B20").Value = "4 channels, 5 comparators: 2376 networks"

# --- Formula tweak on Generator1: B13 now scaled by 0.305, like the other
#     sheets' "1 network (...)" rows -----------------------------------------
$wsGen1.Range("B13").Formula = "=55^33 *0.305"

# --- Highlight the "1 network (...)" result cells with the same yellow
#     fill already used elsewhere in the workbook ---------------------------
$wsGen1.Range("B13").Interior.Color = 65535
$wsGen2.Range("B14").Interior.Color = 65535
$wsGen2.Range("B31").Interior.Color = 65535

# --- Page layout tweaks ------------------------------------------------------
$wsTester.PageSetup.Orientation = 2   # xlLandscape
$wsGen2.PageSetup.Orientation = 2     # xlLandscape
$wsGen1.PageSetup.Orientation = 2     # xlLandscape (newly added pageSetup)

# --- Generator2: narrow the spacer column H and shrink the blank spacer row -
$wsGen2.Columns.Item(8).ColumnWidth = 1.7
$wsGen2.Rows.Item(17).RowHeight = 6

# --- Selections on each sheet (also drives tabSelected / activeTab) --------
$wsTester.Range("B18").Select()
$wsGen1.Range("B2:G13").Select()
$wsGen2.Range("B21").Select()
